$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: Insert two new columns before column D (shifts existing D:K data to F:M) ---
$ws.Range("D:E").EntireColumn.Insert()

# --- Step 2: Copy number/date formatting into the new blank D/E columns from column F/G ---
# (skips label-only rows 5,6,37,79, which have no data columns at all)
$ws.Range("F7:G35").Copy()
$ws.Range("D7:E35").PasteSpecial(-4122)
$ws.Range("F38:G77").Copy()
$ws.Range("D38:E77").PasteSpecial(-4122)
$ws.Range("F80:G102").Copy()
$ws.Range("D80:E102").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Step 3: Write the new D/E quarter values ---
$ws.Range("D7").Value = 43465
$ws.Range("E7").Value = 43373
$ws.Range("D8").Value = 62300
$ws.Range("E8").Value = 58700
$ws.Range("D9").Value = 13900
$ws.Range("E9").Value = 13500
$ws.Range("D10").Value = 48400
$ws.Range("E10").Value = 45200
$ws.Range("D12").Value = 8800
$ws.Range("E12").Value = 7300
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0
$ws.Range("D14").Value = 0
$ws.Range("E14").Value = 0
$ws.Range("D15").Value = 0
$ws.Range("E15").Value = 0
$ws.Range("D17").Value = 70700
$ws.Range("E17").Value = 63800
$ws.Range("D18").Value = -8400
$ws.Range("E18").Value = -5100
$ws.Range("D20").Value = 700
$ws.Range("E20").Value = 600
$ws.Range("D21").Value = -2200
$ws.Range("E21").Value = 1400
$ws.Range("D22").Value = 0
$ws.Range("E22").Value = 0
$ws.Range("D23").Value = -7800
$ws.Range("E23").Value = -4400
$ws.Range("D24").Value = 0
$ws.Range("E24").Value = 0
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 0
$ws.Range("D26").Value = -7800
$ws.Range("E26").Value = -4500
$ws.Range("D27").Value = -7700
$ws.Range("E27").Value = -4500
$ws.Range("D28").Value = 0
$ws.Range("E28").Value = 0
$ws.Range("D29").Value = 0
$ws.Range("E29").Value = 0
$ws.Range("D30").Value = 0
$ws.Range("E30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("E31").Value = 0
$ws.Range("D32").Value = -700
$ws.Range("E32").Value = -600
$ws.Range("D33").Value = -7700
$ws.Range("E33").Value = -4500
$ws.Range("D34").Value = 0
$ws.Range("E34").Value = 0
$ws.Range("D35").Value = -7700
$ws.Range("E35").Value = -4500
$ws.Range("D38").Value = 43465
$ws.Range("E38").Value = 43373
$ws.Range("D41").Value = 46200
$ws.Range("E41").Value = 38400
$ws.Range("D42").Value = 86400
$ws.Range("E42").Value = 86700
$ws.Range("D43").Value = 80100
$ws.Range("E43").Value = 72500
$ws.Range("D44").Value = 0
$ws.Range("E44").Value = 0
$ws.Range("D45").Value = 8800
$ws.Range("E45").Value = 6900
$ws.Range("D46").Value = 221500
$ws.Range("E46").Value = 204500
$ws.Range("D47").Value = 0
$ws.Range("E47").Value = 0
$ws.Range("D48").Value = 13500
$ws.Range("E48").Value = 13200
$ws.Range("D49").Value = 221900
$ws.Range("E49").Value = 224900
$ws.Range("D50").Value = 0
$ws.Range("E50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("E51").Value = 0
$ws.Range("D52").Value = 36900
$ws.Range("E52").Value = 32100
$ws.Range("D53").Value = 0
$ws.Range("E53").Value = 0
$ws.Range("D54").Value = 493900
$ws.Range("E54").Value = 474700
$ws.Range("D57").Value = 3400
$ws.Range("E57").Value = 2900
$ws.Range("D58").Value = 0
$ws.Range("E58").Value = 0
$ws.Range("D59").Value = 155800
$ws.Range("E59").Value = 140400
$ws.Range("D60").Value = 159200
$ws.Range("E60").Value = 143300
$ws.Range("D61").Value = 0
$ws.Range("E61").Value = 0
$ws.Range("D62").Value = 8700
$ws.Range("E62").Value = 8900
$ws.Range("D63").Value = 0
$ws.Range("E63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("E64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("E65").Value = 0
$ws.Range("D66").Value = 172300
$ws.Range("E66").Value = 152300
$ws.Range("D68").Value = 0
$ws.Range("E68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("E69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("E70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("E71").Value = 0
$ws.Range("D72").Value = -130600
$ws.Range("E72").Value = -122900
$ws.Range("D73").Value = 0
$ws.Range("E73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("E74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("E75").Value = 0
$ws.Range("D76").Value = 321600
$ws.Range("E76").Value = 322400
$ws.Range("D77").Value = 0
$ws.Range("E77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("E80").Value = 43373
$ws.Range("D81").Value = -7700
$ws.Range("E81").Value = -4500
$ws.Range("D83").Value = 5600
$ws.Range("E83").Value = 5800
$ws.Range("D84").Value = 0
$ws.Range("E84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("E85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("E86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("E87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("E88").Value = 0
$ws.Range("D89").Value = 4800
$ws.Range("E89").Value = 4800
$ws.Range("D91").Value = -1700
$ws.Range("E91").Value = -900
$ws.Range("D92").Value = 0
$ws.Range("E92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("E93").Value = 0
$ws.Range("D94").Value = -2100
$ws.Range("E94").Value = -2600
$ws.Range("D96").Value = 0
$ws.Range("E96").Value = 0
$ws.Range("D97").Value = 0
$ws.Range("E97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("E98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("E99").Value = 0
$ws.Range("D100").Value = 4800
$ws.Range("E100").Value = 7300
$ws.Range("D101").Value = 0
$ws.Range("E101").Value = 0
$ws.Range("D102").Value = 7800
$ws.Range("E102").Value = 9500

# --- Step 4: Correct the restated historical values that did not purely shift ---
$ws.Range("H8").Value = 50000
$ws.Range("I8").Value = 45400
$ws.Range("H9").Value = 19800
$ws.Range("I9").Value = 19600
$ws.Range("H10").Value = 30200
$ws.Range("I10").Value = 25800
$ws.Range("H17").Value = 54200
$ws.Range("I17").Value = 57700
$ws.Range("H18").Value = -4200
$ws.Range("I18").Value = -12300
$ws.Range("H20").Value = 300
$ws.Range("H21").Value = 1300
$ws.Range("I21").Value = -6900
$ws.Range("H23").Value = -3900
$ws.Range("I23").Value = -12000
$ws.Range("H24").Value = 300
$ws.Range("I24").Value = 100
$ws.Range("H26").Value = -4200
$ws.Range("I26").Value = -12100
$ws.Range("H27").Value = -4200
$ws.Range("I27").Value = -12100
$ws.Range("H32").Value = -300
$ws.Range("H33").Value = -4200
$ws.Range("I33").Value = -12100
$ws.Range("H35").Value = -4200
$ws.Range("I35").Value = -12100
$ws.Range("H81").Value = -4200
$ws.Range("I81").Value = -12100
$ws.Range("F91").Value = -2100
$ws.Range("G91").Value = -1600
$ws.Range("H91").Value = -300
$ws.Range("I91").Value = -2600
$ws.Range("J91").Value = -600
